# Update "想去人数" (want-to-go count) figures in the F column of the
# "展览" and "全部类型" sheets to match the refreshed data snapshot.

$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)

    $ws.Range("F2").Value = 723
    $ws.Range("F3").Value = 4085
    $ws.Range("F4").Value = 115
    $ws.Range("F5").Value = 750
}
